# EEPROM save data now works (with a delay instead of ACK polling),
# added error eeprom report on screen.
#
# This updates the "Typography" sheet of the TouchGFX texts workbook:
#   - Sets the (previously empty) "Ellipsis Character" column (J) to "-"
#     for the Large/Small/XLarge/etc. typography rows (4, 5, 6, 7, 8).
#   - Updates row 6's "Wildcard Characters" (G) to include an underscore,
#     and its "Wildcard Ranges" (I) to also allow A-Z.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Typography")

# Ellipsis Character column (J) for rows 4-8
$ws.Range("J4").Value = "-"
$ws.Range("J5").Value = "-"
$ws.Range("J6").Value = "-"
$ws.Range("J7").Value = "-"
$ws.Range("J8").Value = "-"

# Row 6: Wildcard Characters (G) and Wildcard Ranges (I)
$ws.Range("G6").Value = '.",_'
$ws.Range("I6").Value = "0-9,A-Z"
